$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.027.30'
$ws.Range('E2').Value = '  +1.98%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.705.37'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.59'
$ws.Range('E5').Value = '  -0.16%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3992'
$ws.Range('E7').Value = '  +2.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4045'
$ws.Range('E8').Value = '  -0.79%  '
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.98'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08817'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '26.07'
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.489'
$ws.Range('E14').Value = '  -0.18%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.984'
$ws.Range('E15').Value = '  -3.76%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001351'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.714.14'
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('E18').Value = '  -1.92%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.79'
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('E23').Value = '  -0.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '25.035.09'
$ws.Range('E24').Value = '  +2.00%  '
$ws.Range('E25').Value = '  +2.92%  '
$ws.Range('E26').Value = '  -3.57%  '
$ws.Range('E27').Value = '  +2.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.051'
$ws.Range('E28').Value = '  +12.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.87'
$ws.Range('E29').Value = '  -3.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '152.21'
$ws.Range('E30').Value = '  +2.70%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.448'
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.685'
$ws.Range('E32').Value = '  +21.61%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.08641'
$ws.Range('E33').Value = '  -1.63%  '
$ws.Range('B34').Value = 'VeChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.03166'
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.050'
$ws.Range('E35').Value = '  -0.39%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.206'
$ws.Range('E36').Value = '  -1.39%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2922'
$ws.Range('E37').Value = '  +4.16%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.09714'
$ws.Range('E38').Value = '  +5.61%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.04'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8268'
$ws.Range('E40').Value = '  +3.22%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.04'
$ws.Range('E41').Value = '  -1.68%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.482'
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.03'
$ws.Range('E43').Value = '  -2.70%  '
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.691'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7377'
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.09216'
$ws.Range('E46').Value = '  +12.72%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.250'
$ws.Range('E47').Value = '  -0.58%  '
$ws.Range('B48').Value = 'Flow'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.404'
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('B49').Value = 'Frax'
$ws.Range('C49').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.001'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '140.08'
$ws.Range('E50').Value = '  -0.86%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '92.10'
$ws.Range('E51').Value = '  +1.25%  '
